# Update NATMI LR-pair sheet with newly recalculated TPM values and a new
# "ECs" sending-cluster row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 2 ("Inflammatory-Mac") to make
# room for the new "ECs" sending-cluster entry; existing rows 2-5 shift
# down to rows 3-6.
$ws.Rows.Item(2).Insert()
# The insert picks up formatting from the surrounding rows (bold/centered
# header style); reset the new row back to the plain data-row formatting
# used throughout the rest of the table.
$ws.Range("A2:T2").ClearFormats()

# Row 2: ECs / Ccl3 / Ackr2 / FAPs
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Ccl3"
$ws.Range("C2").Value = "Ackr2"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 1.519357
$ws.Range("H2").Value = 4.558071
$ws.Range("I2").Value = [double]"0.0004273171801484077"
$ws.Range("J2").Value = [double]"0.0004273171801484077"
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.26202
$ws.Range("N2").Value = 0.78606
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.39810192114
$ws.Range("R2").Value = 3.58291729026
$ws.Range("S2").Value = [double]"0.0004273171801484077"
$ws.Range("T2").Value = [double]"0.0004273171801484077"

# Row 3: Inflammatory-Mac / Ccl3 / Ackr2 / FAPs
$ws.Range("A3").Value = "Inflammatory-Mac"
$ws.Range("B3").Value = "Ccl3"
$ws.Range("C3").Value = "Ackr2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 151.6315156666667
$ws.Range("H3").Value = 454.894547
$ws.Range("I3").Value = [double]"0.0426461665667181"
$ws.Range("J3").Value = [double]"0.04264616656671809"
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.26202
$ws.Range("N3").Value = 0.78606
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 39.73048973498
$ws.Range("R3").Value = 357.57440761482
$ws.Range("S3").Value = [double]"0.0426461665667181"
$ws.Range("T3").Value = [double]"0.04264616656671809"

# Row 4: MuSCs / Ccl3 / Ackr2 / FAPs
$ws.Range("A4").Value = "MuSCs"
$ws.Range("B4").Value = "Ccl3"
$ws.Range("C4").Value = "Ackr2"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.08050133333333333
$ws.Range("H4").Value = 0.241504
$ws.Range("I4").Value = [double]"2.264089529859475E-05"
$ws.Range("J4").Value = [double]"2.264089529859474E-05"
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.26202
$ws.Range("N4").Value = 0.78606
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 1
$ws.Range("Q4").Value = 0.02109295936
$ws.Range("R4").Value = 0.18983663424
$ws.Range("S4").Value = [double]"2.264089529859475E-05"
$ws.Range("T4").Value = [double]"2.264089529859474E-05"

# Row 5: Neutrophils / Ccl3 / Ackr2 / FAPs
$ws.Range("A5").Value = "Neutrophils"
$ws.Range("B5").Value = "Ccl3"
$ws.Range("C5").Value = "Ackr2"
$ws.Range("D5").Value = "FAPs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 3259.30721
$ws.Range("H5").Value = 9777.921630000001
$ws.Range("I5").Value = 0.9166759137020294
$ws.Range("J5").Value = 0.9166759137020293
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.26202
$ws.Range("N5").Value = 0.78606
$ws.Range("O5").Value = 1
$ws.Range("P5").Value = 1
$ws.Range("Q5").Value = 854.0036751642
$ws.Range("R5").Value = 7686.0330764778
$ws.Range("S5").Value = 0.9166759137020294
$ws.Range("T5").Value = 0.9166759137020293

# Row 6 (new): Resolving-Mac / Ccl3 / Ackr2 / FAPs
$ws.Range("A6").Value = "Resolving-Mac"
$ws.Range("B6").Value = "Ccl3"
$ws.Range("C6").Value = "Ackr2"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 143.033414
$ws.Range("H6").Value = 429.100242
$ws.Range("I6").Value = 0.04022796165580557
$ws.Range("J6").Value = 0.04022796165580557
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.26202
$ws.Range("N6").Value = 0.78606
$ws.Range("O6").Value = 1
$ws.Range("P6").Value = 1
$ws.Range("Q6").Value = 37.47761513628
$ws.Range("R6").Value = 337.29853622652
$ws.Range("S6").Value = 0.04022796165580557
$ws.Range("T6").Value = 0.04022796165580557
